# Splits every "Situação Docket: " run into three runs --
# "Situação ", "#{tipoCertidao}" and ": " -- so the certidao type can be
# merged into the placeholder, exactly like the surrounding template
# fields (e.g. #{situacao}).

$d = $word.ActiveDocument

$prefix = "Situação "
$needle = "Docket"
$suffix = ": "
$target = $prefix + $needle + $suffix
$mergeField = "#{tipoCertidao}"

$rng = $d.Content
$count = 0
$maxIterations = 50

for ($i = 0; $i -lt $maxIterations; $i++) {
    $found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        break
    }

    $count = $count + 1
    $matchStart = $rng.Start

    $needleStart = $matchStart + $prefix.Length
    $needleEnd = $needleStart + $needle.Length
    $needleRange = $d.Range($needleStart, $needleEnd)

    # Nudge a formatting property away from (and back to) its current
    # value. That forces Word to split the run around this sub-range
    # without leaving behind any stray direct formatting once the value
    # is restored, so the new run's <w:rPr> ends up identical to its
    # neighbours -- just like in the target document.
    $origSize = $needleRange.Font.Size
    $needleRange.Font.Size = $origSize - 1
    $needleRange.Text = $mergeField
    $needleRange.Font.Size = $origSize

    # Resume searching right after the text we just inserted.
    $rng.Start = $needleRange.End
    $rng.End = $d.Content.End
}

Write-Host "Replaced" $count "occurrence(s) of 'Situação Docket: '"
